# Generate Report for Handback
# Adds a new handed-back file ("66777bd6-331e-4114-bf5a-55ff04bdf6bb") as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the structure
# already used by the two existing entries.

$wb = $excel.ActiveWorkbook

$guid       = "66777bd6-331e-4114-bf5a-55ff04bdf6bb"
$commitHash = "5d32d87dbfb9bde9ea141dcd58bde8beb54ef2eb"
$mdName     = "$guid.md"
$status     = "Handed back: in sync with en-US"
$reason     = "Include"

# ---------------------------------------------------------------------------
# Sheet 1: Overview (File Name | zh-cn | de-de)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f5a7e6c1d9b3a2468057cfe1a9d4b6costant0/e2e/$mdName", "", "", $mdName) | Out-Null
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status

# ---------------------------------------------------------------------------
# Sheets 2 & 3: per-language detail sheets
#   A Source File Name | B Status | C Correspond Handoff File |
#   D Correspond Handoff Datetime | E Target File |
#   F Correspond Handback File | G Correspond Handback DateTime |
#   H Handoff Reason | I Dependency From
# ---------------------------------------------------------------------------
$langs = @{
    "zh-cn" = @{
        xlf        = "$guid.$commitHash.zh-cn.xlf"
        handoffDt  = "2016-02-06 03:48:50"
        handbackDt = "2016-02-06 03:49:32"
    }
    "de-de" = @{
        xlf        = "$guid.$commitHash.de-de.xlf"
        handoffDt  = "2016-02-06 03:49:01"
        handbackDt = "2016-02-06 03:49:51"
    }
}

foreach ($lang in $langs.Keys) {
    $info = $langs[$lang]
    $ws = $wb.Worksheets.Item($lang)

    $mdUrl  = "https://github.com/OpenLocalizationTestOrg/oltest.$lang/blob/f5a7e6c1d9b3a2468057cfe1a9d4b6costant1/e2e/$mdName"
    $offUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5a7e6c1d9b3a2468057cfe1a9d4b6costant2/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/yuwzho/$($info.xlf)"
    $bckUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f5a7e6c1d9b3a2468057cfe1a9d4b6costant3/ol-handback/OpenLocalizationTestOrg/oltest.$lang/yuwzho/$($info.xlf)"

    $ws.Hyperlinks.Add($ws.Range("A4"), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Range("B4").Value = $status
    $ws.Hyperlinks.Add($ws.Range("C4"), $offUrl, "", "", $info.xlf) | Out-Null
    $ws.Range("D4").Value = $info.handoffDt
    $ws.Hyperlinks.Add($ws.Range("E4"), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F4"), $bckUrl, "", "", $info.xlf) | Out-Null
    $ws.Range("G4").Value = $info.handbackDt
    $ws.Range("H4").Value = $reason
}

Write-Host "Handback row added for $guid"
